$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B4").Copy()
$ws.Range("B5").PasteSpecial(-4122)  # xlPasteFormats
Write-Host "done"
